# Insert a new column A ("BaseReportCriteriaId") before the existing
# Category / FriendlyField / Field table, shifting the current A:C data
# to B:D, then populate the new column with a header + sequential ids.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift existing columns A:C to B:D by inserting a new column at A.
$ws.Range("A1").EntireColumn.Insert()

# Header for the new id column.
$ws.Range("A1").Value = "BaseReportCriteriaId"

# Sequential ids 1..40 for data rows 2..41.
for ($i = 2; $i -le 41; $i++) {
    $ws.Cells.Item($i, 1).Value = $i - 1
}

# Match the saved view state: scrolled so row 25 is at the top, with the
# new id column (A2:A41) selected.
$ws.Range("A2:A41").Select()
$excel.ActiveWindow.ScrollRow = 25
